# Weekly price update: insert a new data row at row 71, shifting the
# existing rows (old 71..159) down to (72..160), then fill the new
# row 71 with the latest week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 71 - this pushes every
# existing row from 71 downward by one position.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with this week's values. Columns
# that repeat the same constant for every record (A, B, C, E, F, G, H,
# N, Q, R) are copied from the neighboring row; the rest come from the
# new data point.
$ws.Range("A71").Value = 9
$ws.Range("B71").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C71").Value = "Metropolitana"
$ws.Range("D71").Value = 44494
$ws.Range("E71").Value = 13
$ws.Range("F71").Value = 100112026
$ws.Range("G71").Value = "Haba"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 120
$ws.Range("K71").Value = 6000
$ws.Range("L71").Value = 6000
$ws.Range("M71").Value = 6000
$ws.Range("N71").Value = "$/saco 25 kilos"
$ws.Range("O71").Value = "Región Metropolitana"
$ws.Range("P71").Value = 240
$ws.Range("Q71").Value = 25
$ws.Range("R71").Value = "Hortaliza"
